$wb = $excel.ActiveWorkbook

# Trade #6 closed at 2026-02-16 22:56:57 - base_strategy UP +0.000%
# Same new row is appended to both the "All Trades" log and the
# per-strategy "base_strategy" log.
foreach ($sheetName in @("All Trades", "base_strategy")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(7, 1).Value = 6

    # Date/Time columns look like dates/times and would otherwise be
    # auto-converted to date serial numbers - enter them as literal text
    # (leading apostrophe, same as typing '2026-02-16 into Excel) and then
    # drop back to the Normal cell style so no stray number format sticks.
    $ws.Cells.Item(7, 2).Value = "'2026-02-16"
    $ws.Cells.Item(7, 2).Style = "Normal"
    $ws.Cells.Item(7, 3).Value = "'22:56:57"
    $ws.Cells.Item(7, 3).Style = "Normal"

    $ws.Cells.Item(7, 4).Value = "base_strategy"
    $ws.Cells.Item(7, 5).Value = "UP"
    $ws.Cells.Item(7, 6).Value = 0.5

    # Exit Price is blank - the trade is still OPEN. Force an explicit
    # empty text cell (matches the existing blank cells in this column)
    # rather than leaving the cell completely absent.
    $ws.Cells.Item(7, 7).Value = "'"
    $ws.Cells.Item(7, 7).Style = "Normal"

    $ws.Cells.Item(7, 8).Value = "OPEN"
    $ws.Cells.Item(7, 9).Value = 0
    $ws.Cells.Item(7, 10).Value = 0
    $ws.Cells.Item(7, 11).Value = 100
    $ws.Cells.Item(7, 12).Value = 0
    $ws.Cells.Item(7, 13).Value = 0
    $ws.Cells.Item(7, 14).Value = 0.6
    $ws.Cells.Item(7, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason is blank too (trade still OPEN).
    $ws.Cells.Item(7, 16).Value = "'"
    $ws.Cells.Item(7, 16).Style = "Normal"

    $ws.Cells.Item(7, 17).Value = 0
}
